$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (13) for the book currently being read.
$ws.Range("A13").Value = "The Passionate Programmer: Creating a Remarkable Career in Software Development"
$ws.Range("B13").Value = "Chad Fowler"
$ws.Range("C13").Value = 89

# Update the selection to match the author's final cursor position.
$ws.Range("G13").Select()
